$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / header text updates ---
$ws.Range("A8").Value = "Volume 32   Number  26"
$ws.Range("C9").Value = "Report Covering the Week  6/23/2025  Through  6/29/2025"

# --- Weekly crime statistics updates ---
# Row 15
$ws.Range("C14").Copy($ws.Range("D15"))
$ws.Range("H14").Copy($ws.Range("E15"))
$ws.Range("L15").Value = -20
$ws.Range("M15").Value = 14.285714285714

# Row 16
$ws.Range("C16").Value = 5
$ws.Range("E16").Value = 150
$ws.Range("F16").Value = 12
$ws.Range("H16").Value = 140
$ws.Range("I16").Value = 58
$ws.Range("J16").Value = 41
$ws.Range("K16").Value = 41.463414634146
$ws.Range("L16").Value = -6.451612903225
$ws.Range("M16").Value = -56.716417910447
$ws.Range("N16").Value = -92.944038929440

# Row 17
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 250
$ws.Range("F17").Value = 28
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = 55.555555555555
$ws.Range("I17").Value = 180
$ws.Range("J17").Value = 130
$ws.Range("K17").Value = 38.461538461538
$ws.Range("L17").Value = 11.111111111111
$ws.Range("M17").Value = 24.137931034482
$ws.Range("N17").Value = -60.176991150442

# Row 18
$ws.Range("G15").Copy($ws.Range("C18"))
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 100
$ws.Range("I18").Value = 35
$ws.Range("J18").Value = 38
$ws.Range("K18").Value = -7.894736842105
$ws.Range("L18").Value = -2.777777777777
$ws.Range("M18").Value = -75.177304964539
$ws.Range("N18").Value = -96.268656716417

# Row 19
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = -33.333333333333
$ws.Range("F19").Value = 19
$ws.Range("G19").Value = 18
$ws.Range("H19").Value = 5.555555555555
$ws.Range("I19").Value = 149
$ws.Range("J19").Value = 122
$ws.Range("K19").Value = 22.131147540983
$ws.Range("L19").Value = -31.018518518518
$ws.Range("M19").Value = -24.747474747474
$ws.Range("N19").Value = -51.307189542483

# Row 20
$ws.Range("J14").Copy($ws.Range("C20"))
$ws.Range("C14").Copy($ws.Range("D20"))
$ws.Range("H14").Copy($ws.Range("E20"))
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 10
$ws.Range("I20").Value = 44
$ws.Range("K20").Value = 10
$ws.Range("L20").Value = -12
$ws.Range("M20").Value = -45
$ws.Range("N20").Value = -93.432835820895

# Row 21
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 11
$ws.Range("E21").Value = 100
$ws.Range("F21").Value = 77
$ws.Range("G21").Value = 62
$ws.Range("H21").Value = 24.193548387096
$ws.Range("I21").Value = 475
$ws.Range("J21").Value = 388
$ws.Range("K21").Value = 22.422680412371
$ws.Range("L21").Value = -11.873840445269
$ws.Range("M21").Value = -32.814710042432
$ws.Range("N21").Value = -85.362095531587

# Row 22
$ws.Range("C14").Copy($ws.Range("D22"))
$ws.Range("H14").Copy($ws.Range("E22"))
$ws.Range("M22").Value = -96

# Row 24
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 11
$ws.Range("E24").Value = 54.545454545454
$ws.Range("F24").Value = 75
$ws.Range("G24").Value = 48
$ws.Range("H24").Value = 56.25
$ws.Range("I24").Value = 480
$ws.Range("J24").Value = 423
$ws.Range("K24").Value = 13.475177304964
$ws.Range("L24").Value = -12.087912087912
$ws.Range("M24").Value = 2.783725910064

# Row 25
$ws.Range("C25").Value = 4
$ws.Range("C14").Copy($ws.Range("D25"))
$ws.Range("H14").Copy($ws.Range("E25"))
$ws.Range("F25").Value = 14
$ws.Range("G25").Value = 5
$ws.Range("H25").Value = 180
$ws.Range("I25").Value = 129
$ws.Range("K25").Value = 72
$ws.Range("L25").Value = 11.206896551724

# Row 26
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 12
$ws.Range("E26").Value = -8.333333333333
$ws.Range("F26").Value = 50
$ws.Range("G26").Value = 44
$ws.Range("H26").Value = 13.636363636363
$ws.Range("I26").Value = 243
$ws.Range("J26").Value = 266
$ws.Range("K26").Value = -8.646616541353
$ws.Range("L26").Value = -4.330708661417
$ws.Range("M26").Value = -39.097744360902

# Row 27
$ws.Range("C14").Copy($ws.Range("D27"))
$ws.Range("H14").Copy($ws.Range("E27"))
$ws.Range("L27").Value = -33.333333333333

# Row 28
$ws.Range("C14").Copy($ws.Range("C28"))
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = -100
$ws.Range("G28").Value = 12
$ws.Range("H28").Value = -75
$ws.Range("J28").Value = 39
$ws.Range("K28").Value = -28.205128205128
$ws.Range("L28").Value = 12

# Row 29
$ws.Range("M29").Value = -82.352941176470
$ws.Range("N29").Value = -95

# Row 30
$ws.Range("M30").Value = -76.923076923076
$ws.Range("N30").Value = -94.915254237288

# Row 31
$ws.Range("I14").Copy($ws.Range("F31"))
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 9
$ws.Range("K31").Value = 80
$ws.Range("L31").Value = 50
